# Updates the cryptos list with refreshed prices / 1h volume percentages,
# matching the "Updated cryptos list ... with GitHub Actions" commit.
# Note: some Price (column D) values look like plain decimals (e.g. "229.90")
# and Excel would otherwise auto-convert them to numbers; a leading apostrophe
# forces them to stay as text, matching the original cell content type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '93.322.72'
$ws.Range("E2").Value = '  -3.16%  '
$ws.Range("D3").Value = '3.318.72'
$ws.Range("E3").Value = '  -4.95%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '''229.90'
$ws.Range("E5").Value = '  -6.17%  '
$ws.Range("D6").Value = '''618.12'
$ws.Range("E6").Value = '  -4.91%  '
$ws.Range("E7").Value = '  -4.13%  '
$ws.Range("D8").Value = '''0.384'
$ws.Range("E8").Value = '  -7.53%  '
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("D10").Value = '''0.920'
$ws.Range("E10").Value = '  -8.75%  '
$ws.Range("D11").Value = '3.318.21'
$ws.Range("E11").Value = '  -4.84%  '
$ws.Range("D12").Value = '''41.53'
$ws.Range("E12").Value = '  -4.34%  '
$ws.Range("E13").Value = '  -4.10%  '
$ws.Range("E14").Value = '  -3.50%  '
$ws.Range("D15").Value = '93.128.15'
$ws.Range("E15").Value = '  -3.04%  '
$ws.Range("D16").Value = '3.942.37'
$ws.Range("E16").Value = '  -4.72%  '
$ws.Range("D17").Value = '''0.0000242'
$ws.Range("E17").Value = '  -5.20%  '
$ws.Range("D18").Value = '''7.96'
$ws.Range("E18").Value = '  -7.68%  '
$ws.Range("D19").Value = '3.319.83'
$ws.Range("E19").Value = '  -4.82%  '
$ws.Range("D20").Value = '''17.12'
$ws.Range("E20").Value = '  -7.42%  '
$ws.Range("D21").Value = '''10.86'
$ws.Range("E21").Value = '  -10.87%  '
$ws.Range("D22").Value = '''3.40'
$ws.Range("E22").Value = '  +2.47%  '
$ws.Range("D23").Value = '''490.52'
$ws.Range("E23").Value = '  -4.98%  '
$ws.Range("D24").Value = '''0.450'
$ws.Range("E24").Value = '  -9.97%  '
$ws.Range("E25").Value = '  -7.71%  '
$ws.Range("D26").Value = '''6.02'
$ws.Range("E26").Value = '  -10.65%  '
$ws.Range("D27").Value = '''89.48'
$ws.Range("E27").Value = '  -3.20%  '
$ws.Range("D28").Value = '3.499.80'
$ws.Range("E28").Value = '  -4.44%  '
$ws.Range("D29").Value = '''11.53'
$ws.Range("E29").Value = '  -8.33%  '
$ws.Range("E30").Value = '  +0.16%  '
$ws.Range("D31").Value = '''10.98'
$ws.Range("E31").Value = '  -8.84%  '
$ws.Range("E32").Value = '  -3.86%  '
$ws.Range("E33").Value = '  -6.17%  '
$ws.Range("E34").Value = '  +0.12%  '
$ws.Range("E35").Value = '  -7.70%  '
$ws.Range("D36").Value = '''28.22'
$ws.Range("E36").Value = '  -10.18%  '
$ws.Range("D37").Value = '''0.526'
$ws.Range("E37").Value = '  -10.21%  '
$ws.Range("D38").Value = '''521.76'
$ws.Range("E38").Value = '  +0.15%  '
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("E40").Value = '  -8.08%  '
$ws.Range("E41").Value = '  -4.38%  '
$ws.Range("E42").Value = '  -9.47%  '
$ws.Range("E43").Value = '  -8.12%  '
$ws.Range("E44").Value = '  -0.65%  '
$ws.Range("D45").Value = '''1.66'
$ws.Range("E45").Value = '  -4.93%  '

# Rows 46/47 swapped rank order: VeChain now ranks above MantraDAO.
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '''0.0405'
$ws.Range("E46").Value = '  -4.59%  '
$ws.Range("B47").Value = 'MantraDAO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D47").Value = '''3.55'
$ws.Range("E47").Value = '  -2.59%  '
$ws.Range("D48").Value = '''5.31'
$ws.Range("E48").Value = '  -5.26%  '
$ws.Range("D49").Value = '''52.65'
$ws.Range("E49").Value = '  -0.63%  '
$ws.Range("D50").Value = '''2.09'
$ws.Range("E50").Value = '  -5.08%  '
$ws.Range("D51").Value = '''7.85'
$ws.Range("E51").Value = '  -5.76%  '
